# Auto-generated Word COM-interop script to apply the scripture-of-the-day update.
# The document body is a single paragraph / single run consisting of alternating
# <w:t>verse text</w:t><w:br/> pairs. The diff replaces every existing verse line
# with a new one (same position) and appends 25 brand-new verse lines at the end.
$d = $word.ActiveDocument

# Sequentially replace each dated scripture line with its new-date counterpart.
# Matches are resolved strictly in document order: after each replacement the next
# search only looks in the remainder of the document (from $pos onward). This keeps
# short / ambiguous text such as a lone verse number (e.g. "32.") from matching a
# later, unrelated verse that merely happens to start with the same digits.
$pos = 0
$pairs = @(
    ,@("********************************8月31日读经章节***************************", "********************************9月1日读经章节***************************")
    ,@("Chapter 4 of 1_Thessalonians", "Chapter 5 of 1_Thessalonians")
    ,@("1.弟兄们，我还有话说。我们靠着主耶稣求你们，劝你们，你们既然受了我们的教训，知道该怎样行，可以讨神的喜悦，就要照你们现在所行的，更加勉励。", "1.弟兄们，论到时候日期，不用写信给你们。")
    ,@("2.你们原晓得我们凭主耶稣传给你们什么命令。", "2.因为你们自己明明晓得，主的日子来到，好像夜间的贼一样。")
    ,@("3.神的旨意就是要你们成为圣洁，远避淫行。", "3.人正说平安稳妥的时候，灾祸忽然临到他们，如同产难临到怀胎的妇人一样，他们绝不能逃脱。")
    ,@("4.要你们各人晓得怎样用圣洁尊贵，守着自己的身体。", "4.弟兄们，你们却不在黑暗里，叫那日子临到你们像贼一样。")
    ,@("5.不放纵私欲的邪情，像那不认识神的外邦人。", "5.你们都是光明之子，都是白昼之子，我们不是属黑夜的，也不是属幽暗的。")
    ,@("6.不要一个人在这事上越分，欺负他的弟兄。因为这一类的事，主必报应，正如我预先对你们说过，又切切嘱咐你们的。", "6.所以我们不要睡觉，像别人一样，总要儆醒谨守。")
    ,@("7.神召我们，本不是要我们沾染污秽，乃是要我们成为圣洁。", "7.因为睡了的人是在夜间睡。醉了的人是在夜间醉。")
    ,@("8.所以那弃绝的，不是弃绝人，乃是弃绝那赐圣灵给你们的神。", "8.但我们既然属乎白昼，就应当谨守，把信和爱当作护心镜遮胸。把得救的盼望当作头盔戴上。")
    ,@("9.论到弟兄们相爱，不用人写信给你们，因为你们自己蒙了神的教训，叫你们彼此相爱。", "9.因为神不是预定我们受刑，乃是预定我们借着我们主耶稣基督得救。")
    ,@("10.你们向马其顿全地的众弟兄，固然是这样行，但我劝弟兄们要更加勉励。", "10.他替我们死，叫我们无论醒着睡着，都与他同活。")
    ,@("11.又要立志作安静人，办自己的事，亲手作工，正如我们从前所吩咐你们的。", "11.所以你们该彼此劝慰，互相建立，正如你们素常所行的。")
    ,@("12.叫你们可以向外人行事端正，自己也就没有什么缺乏了。", "12.弟兄们，我们劝你们敬重那在你们中间劳苦的人，就是在主里面治理你们，劝戒你们的。")
    ,@("13.论到睡了的人，我们不愿意弟兄们不知道，恐怕你们忧伤，像那些没有指望的人一样。", "13.又因他们所作的工，用爱心格外尊重他们，你们也要彼此和睦。")
    ,@("14.我们若信耶稣死而复活了，那已经在耶稣里睡了的人，神也必将他与耶稣一同带来。", "14.我们又劝弟兄们，要警戒不守规矩的人。勉励灰心的人。扶助软弱的人。也要向众人忍耐。")
    ,@("15.我们现在照主的话告诉你们一件事。我们这活着还存留到主降临的人，断不能在那已经睡了的人之先。", "15.你们要谨慎，无论是谁都不可以恶报恶。或是彼此相待，或是待众人，常要追求良善。")
    ,@("16.因为主必亲自从天降临，有呼叫的声音，和天使长的声音，又有神的号吹响。那在基督里死了的人必先复活。", "16.要常常喜乐。")
    ,@("17.以后我们这活着还存留的人，必和他们一同被提到云里，在空中与主相遇。这样，我们就要和主永远同在。", "17.不住地祷告。")
    ,@("18.所以你们当用这些话彼此劝慰。", "18.凡事谢恩。因为这是神在基督耶稣里向你们所定的旨意。")
    ,@("Chapter 23 of Proverbs", "19.不要销灭圣灵的感动。")
    ,@("1.你若与官长坐席，要留意在你面前的是谁。", "20.不要藐视先知的讲论。")
    ,@("2.你若是贪食的，就当拿刀放在喉咙上。", "21.但要凡事察验。善美的要持守。")
    ,@("3.不可贪恋他的美食，因为是哄人的食物。", "22.各样的恶事要禁戒不作。")
    ,@("4.不要劳碌求富。休仗自己的聪明。", "23.愿赐平安的神，亲自使你们全然成圣。又愿你们的灵，与魂，与身子，得蒙保守，在我主耶稣基督降临的时候，完全无可指摘。")
    ,@("5.你岂要定睛在虚无的钱财上吗？因钱财必长翅膀，如鹰向天飞去。", "24.那召你们的本是信实的，他必成就这事。")
    ,@("6.不要吃恶眼人的饭。也不要贪他的美味。", "25.请弟兄们为我们祷告。")
    ,@("7.因为他心怎样思量，他为人就是怎样。他虽对你说，请吃，请喝。他的心却与你相背。", "26.与众弟兄亲嘴问安务要圣洁。")
    ,@("8.你所吃的那点食物，必吐出来。你所说的甘美言语，也必落空。", "27.我指着主嘱咐你们，要把这信念给众弟兄听。")
    ,@("9.你不要说话给愚昧人听。因他必藐视你智慧的言语。", "28.愿我主耶稣基督的恩常与你们同在。")
    ,@("10.不可挪移古时的地界。也不可侵入孤儿的田地。", "Chapter 25 of Proverbs")
    ,@("11.因他们的救赎主，大有能力。他必向你为他们辨屈。", "1.以下也是所罗门的箴言。是犹大王希西家的人所誊录的。")
    ,@("12.你要留心领受训诲。侧耳听从知识的言语。", "2.将事隐秘，乃神的荣耀。将事察清，乃君王的荣耀。")
    ,@("13.不可不管教孩童，你用杖打他，他必不至于死。", "3.天之高，地之厚，君王之心也测不透。")
    ,@("14.你要用杖打他，就可以救他的灵魂免下阴间。", "4.除去银子的渣滓，就有银子出来，银匠能以作器皿。")
    ,@("15.我儿你心若存智慧，我的心也甚欢喜。", "5.除去王面前的恶人，国位就靠公义坚立。")
    ,@("16.你的嘴若说正直话，我的心肠也必快乐。", "6.不要在王面前妄自尊大。不要在大人的位上站立。")
    ,@("17.你心中不要嫉妒罪人。只要终日敬畏耶和华。", "7.宁可有人说，请你上来，强如在你觐见的王子面前，叫你退下。")
    ,@("18.因为至终必有善报。你的指望也不至断绝。", "8.不要冒失出去与人争竞，免得至终被他羞辱，你就不知道怎样行了。")
    ,@("19.我儿，你当听，当存智慧，好在正道上引导你的心。", "9.你与邻舍争讼，要与他一人辩论。不可泄漏人的密事。")
    ,@("20.好饮酒的，好吃肉的，不要与他们来往。", "10.恐怕听见的人骂你，你的臭名就难以脱离。")
    ,@("21.因为好酒贪食的，必致贫穷。好睡觉的，必穿破烂衣服。", "11.一句话说得合宜，就如金苹果在银网子里。")
    ,@("22.你要听从生你的父亲。你母亲老了，也不可藐视她。", "12.智慧人的劝戒，在顺从的人耳中，好像金耳环，和精金的妆饰。")
    ,@("23.你当买真理。就是智慧，训诲，和聪明，也都不可卖。", "13.忠信的使者，叫差他的人心里舒畅，就如在收割时，有冰雪的凉气。")
    ,@("24.义人的父亲，必大得快乐。人生智慧的儿子，必因他欢喜。", "14.空夸赠送礼物的，好像无雨的风云。")
    ,@("25.你要使父母欢喜。使生你的快乐。", "15.恒常忍耐，可以劝动君王。柔和的舌头，能折断骨头。")
    ,@("26.我儿，要将你的心归我。你的眼目，也要喜悦我的道路。", "16.你得了蜜吗？只可吃够而已。恐怕你过饱就呕吐出来。")
    ,@("27.妓女是深坑。外女是窄阱。", "17.你的脚要少进邻舍的家，恐怕他厌烦你，恨恶你。")
    ,@("28.她埋伏好像强盗，她使人中多有奸诈的。", "18.作假见证陷害邻舍的，就是大槌，是利刀，是快箭。")
    ,@("29.谁有祸患，谁有忧愁，谁有争斗，谁有哀叹，（或作怨言）谁无故受伤，谁眼目红赤。", "19.患难时倚靠不忠诚的人，好像破坏的牙，错骨缝的脚。")
    ,@("30.就是那流连饮酒，常去寻找调和酒的人。", "20.对伤心的人唱歌，就如冷天脱衣服，又如碱上倒醋。")
    ,@("31.酒发红，在杯中闪烁，你不可观看，虽然下咽舒畅，终久是咬你如蛇，刺你如毒蛇。", "21.你的仇敌，若饿了就给他饭吃。若渴了就给他水喝。")
    ,@("32.", "22.因为你这样行，就是把炭火堆在他的头上。耶和华也必赏赐你。")
    ,@("33.你眼必看见异怪的事。（异怪的事或作淫妇）你心必发出乖谬的话。", "23.北风生雨，谗谤人的舌头也生怒容。")
    ,@("34.你必像躺在海中，或像卧在桅杆上。", "24.宁可住在房顶的角上，不在宽阔的房屋与争吵的妇人同住。")
    ,@("35.你必说，人打我，我却未受伤，人鞭打我，我竟不觉得，我几时清醒，我仍去寻酒。", "25.有好消息从远方来，就如拿凉水给口渴的人喝。")
    ,@("Chapter 24 of Proverbs", "26.义人在恶人面前退缩，好像？？浑之泉，弄浊之井。")
    ,@("1.你不要嫉妒恶人，也不要起意与他们相处。", "27.吃蜜过多，是不好的。考究自己的荣耀，也是可厌的。")
    ,@("2.因为他们的心，图谋强暴。他们的口谈论奸恶。", "28.人不制伏自己的心，好像毁坏的城邑，没有墙垣。")
    ,@("3.房屋因智慧建造，又因聪明立稳。", "Chapter 26 of Proverbs")
    ,@("4.其中因知识充满各样美好宝贵的财物。", "1.夏天落雪，收割时下雨，都不相宜。愚昧人得尊荣，也是如此。")
    ,@("5.智慧人大有能力。有知识的人，力上加力。", "2.麻雀往来，燕子翻飞，这样，无故地咒诅，也必不临到。")
    ,@("6.你去打仗，要凭智谋。谋士众多，人便得胜。", "3.鞭子是为打马。辔头是为勒驴。刑杖是为打愚昧人的背。")
    ,@("7.智慧极高，非愚昧人所能及，所以在城门内，不敢开口。", "4.不要照愚昧人的愚妄话回答他，恐怕你与他一样。")
    ,@("8.设计作恶的，必称为奸人。", "5.要照愚昧人的愚妄话回答他，免得他自以为有智慧。")
    ,@("9.愚妄人的思念，乃是罪恶。亵慢者为人所憎恶。", "6.借愚昧人手寄信的，是砍断自己的脚，自受损害。（自受原文作喝）")
    ,@("10.你在患难之日若胆怯，你的力量就微小。", "7.瘸子的脚，空存无用。箴言在愚昧人的口中，也是如此。")
    ,@("11.人被拉到死地，你要解救。人将被杀，你须拦阻。", "8.将尊荣给愚昧人的，好像人把石子包在机弦里。")
    ,@("12.你若说，这事我未曾知道。那衡量人心的，岂不明白吗？保守你命的，岂不知道吗？他岂不按各人所行的，报应各人吗？", "9.箴言在愚昧人的口中，好像荆棘刺入醉汉的手。")
    ,@("13.我儿，你要吃蜜，因为是好的。吃蜂房下滴的蜜，便觉甘甜。", "10.雇愚昧人的，与雇过路人的，就像射伤众人的弓箭手。")
    ,@("14.你心得了智慧，也必觉得如此。你若找着，至终必有善报。你的指望，也不至断绝。", "11.愚昧人行愚妄事，行了又行，就如狗转过来吃它所吐的。")
    ,@("15.你这恶人，不要埋伏攻击义人的家。不要毁坏他安居之所。", "12.你见自以为有智慧的人吗？愚昧人比他更有指望。")
    ,@("16.因为义人虽七次跌倒，仍必兴起。恶人却被祸患倾倒。", "13.懒惰人说，道上有猛狮，街上有壮狮。")
    ,@("17.你仇敌跌倒，你不要欢喜。他倾倒，你心不要快乐。", "14.门在枢纽转动，懒惰人在床上也是如此。")
    ,@("18.恐怕耶和华看见就不喜悦，将怒气从仇敌身上转过来。", "15.懒惰人放手在盘子里，就是向口撤回，也以为劳乏。")
    ,@("19.不要为作恶的心怀不平。也不要嫉妒恶人。", "16.懒惰人看自己，比七个善于应对的人更有智慧。")
    ,@("20.因为恶人终不得善报。恶人的灯也必熄灭。", "17.过路被事激动，管理不干己的争竞，好像人揪住狗耳。")
    ,@("21.我儿，你要敬畏耶和华与君王。不要与反覆无常的人结交。", "18.人欺凌邻舍，却说，我岂不是戏耍吗？他就像疯狂的人，抛掷火把，利箭，与杀人的兵器。（杀人的兵器原文作死亡）")
    ,@("22.因为他们的灾难，必忽然而起。耶和华与君王所施行的毁灭，谁能知道呢？", "19.")
    ,@("23.以下也是智慧人的箴言。审判时看人情面，是不好的。", "20.火缺了柴，就必熄灭。无人传舌，争竞便止息。")
    ,@("24.对恶人说，你是义人的，这人万民必咒诅，列邦必憎恶。", "21.好争竞的人煽惑争端，就如余火加炭，火上加柴一样。")
    ,@("25.责备恶人的，必得喜悦。美好的福，也必临到他。", "22.传舌人的言语，如同美食，深入人的心腹。")
    ,@("26.应对正直的，犹如与人亲嘴。", "23.火热的嘴，奸恶的心，好像银渣包的瓦器。")
    ,@("27.你要在外头预备工料，在田间办理整齐，然后建造房屋。", "24.怨恨人的用嘴粉饰，心里却藏着诡诈。")
    ,@("28.不可无故作见证，陷害邻舍。也不可用嘴欺骗人。", "25.他用甜言蜜语，你不可信他，因为他心中有七样可憎恶的。")
    ,@("29.不可说，人怎样待我，我也怎样待他，我必照他所行的报复他。", "26.他虽用诡诈遮掩自己的怨恨，他的邪恶必在会中显露。")
    ,@("30.我经过懒惰人的田地，无知人的葡萄园。", "27.挖陷坑的，自己必掉在其中。滚石头的，石头必反滚在他身上。")
    ,@("31.荆棘长满了地皮，刺草遮盖了田面，石墙也坍塌了。", "28.虚谎的舌，恨他所压伤的人。谄媚的口，败坏人的事。")
    ,@("32.我看见就留心思想，我看着就领了训诲。", "Chapter 27 of Proverbs")
    ,@("33.再睡片时，打盹片时，抱着手躺卧片时，", "1.不要为明日自夸，因为一日要生何事，你尚且不能知道。")
    ,@("34.你的贫穷，就必如强盗速来，你的缺乏，仿佛拿兵器的人来到。", "2.要别人夸奖你，不可用口自夸。等外人称赞你，不可用嘴自称。")
)

foreach ($pair in $pairs) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $rng = $d.Range($pos, $d.Content.End)
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    if (-not $found) {
        throw "Could not find expected scripture text: $oldText"
    }
    $pos = $rng.End
}

# Append the newly-added verses (Proverbs 27:3-27) at the very end of the
# document, right after the last existing line and before the paragraph's
# trailing line break, preserving the existing text/line-break pattern.
$newVerses = @(
    ,"3.石头重，沙土沉，愚妄人的恼怒，比这两样更重。"
    ,"4.忿怒为残忍，怒气为狂澜，惟有嫉妒，谁能敌得住呢？"
    ,"5.当面的责备，强如背地的爱情。"
    ,"6.朋友加的伤痕，出于忠诚，仇敌连连亲嘴，却是多余。"
    ,"7.人吃饱了，厌恶蜂房的蜜。人饥饿了，一切苦物都觉甘甜。"
    ,"8.人离本处飘流，好像雀鸟离窝游飞。"
    ,"9.膏油与香料，使人心喜悦。朋友诚实的劝教，也是如此甘美。"
    ,"10.你的朋友，和父亲的朋友，你都不可离弃。你遭难的日子，不要上弟兄的家去。相近的邻舍，强如远方的弟兄。"
    ,"11.我儿，你要作智慧人，好叫我的心欢喜，使我可以回答那讥诮我的人。"
    ,"12.通达人见祸藏躲。愚蒙人前往受害。"
    ,"13.谁为生人作保，就拿谁的衣服。谁为外女作保，谁就承当。"
    ,"14.清晨起来，大声给朋友祝福的，就算是咒诅他。"
    ,"15.大雨之日连连滴漏，和争吵的妇人一样。"
    ,"16.想拦阻她的，便是拦阻风，也是右手抓油。"
    ,"17.铁磨铁，磨出刃来。朋友相感，（原文作磨朋友的脸）也是如此。"
    ,"18.看守无花果树的，必吃树上的果子。敬奉主人的，必得尊荣。"
    ,"19.水中照脸，彼此相符。人与人，心也相对。"
    ,"20.阴间和灭亡，永不满足。人的眼目，也是如此。"
    ,"21.鼎为炼银，炉为炼金，人的称赞也试炼人。"
    ,"22.你虽用杵，将愚妄人与打碎的麦子一同捣在臼中，他的愚妄还是离不了他。"
    ,"23.你要详细知道你羊群的景况。留心料理你的牛群。"
    ,"24.因为赀财不能永有。冠冕岂能存到万代。"
    ,"25.干草割去，嫩草发现，山上的菜蔬，也被收敛。"
    ,"26.羊羔之毛，是为你作衣服。山羊是为作田地的价值。"
    ,"27.并有母山羊奶够你吃，也够你的家眷吃，且够养你的婢女。"
)

$endPos = $d.Content.End - 1
$tail = $d.Range($endPos, $endPos)
foreach ($verse in $newVerses) {
    $tail.InsertAfter($verse + "`v")
    $tail.Collapse(0)
}

